# Updates cryptos list data (Price and Volume(1h) columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.142.99"
$ws.Range("E2").Value = "  -4.83%  "
$ws.Range("D3").Value = "2.998.80"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'570.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.55%  "
$ws.Range("D6").Value = "'125.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.33%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "2.995.42"
$ws.Range("E8").Value = "  -5.17%  "
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("E10").Value = "  -7.23%  "
$ws.Range("D11").Value = "'5.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.74%  "
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("E13").Value = "  -7.44%  "
$ws.Range("D14").Value = "'32.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.90%  "
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "3.500.44"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").Value = "2.996.39"
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("D18").Value = "60.193.31"
$ws.Range("E18").Value = "  -4.72%  "
$ws.Range("D19").Value = "'6.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "'429.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.57%  "
$ws.Range("D21").Value = "'13.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.59%  "
$ws.Range("D22").Value = "'0.672"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("D23").Value = "'7.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.28%  "
$ws.Range("D24").Value = "'12.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").Value = "'79.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.24%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -6.11%  "
$ws.Range("D29").Value = "'1.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.87%  "
$ws.Range("D30").Value = "'7.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.52%  "
$ws.Range("D31").Value = "'6.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.88%  "
$ws.Range("D32").Value = "'25.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.74%  "
$ws.Range("D33").Value = "'0.0952"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.11%  "
$ws.Range("D34").Value = "'5.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("D35").Value = "'0.931"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.66%  "
$ws.Range("D36").Value = "'50.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").Value = "'2.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -15.79%  "
$ws.Range("D38").Value = "'8.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.39%  "
$ws.Range("D39").Value = "0.0₃0657"
$ws.Range("E39").Value = "  -10.32%  "
$ws.Range("E40").Value = "  -8.25%  "
$ws.Range("D41").Value = "'0.107"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("D42").Value = "'372.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("D43").Value = "2.669.29"
$ws.Range("E43").Value = "  -4.27%  "
$ws.Range("E44").Value = "  -6.72%  "
$ws.Range("D46").Value = "'121.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.53%  "
$ws.Range("E47").Value = "  -6.66%  "
$ws.Range("D48").Value = "'1.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").Value = "'23.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.68%  "
$ws.Range("E51").Value = "  -2.10%  "
